$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price column (D) updates - force Text format so numeric-looking strings
# (e.g. "3.00", "70.65") are preserved exactly as text, matching the source data.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '37.705.93'
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.085.86'
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '232.17'
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '57.16'
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0776'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '2.383.29'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '14.40'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '5.21'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '2.097.00'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '37.658.46'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '6.13'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '70.65'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.0₃0819'
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '227.76'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '168.00'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '8.90'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '19.43'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.60'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.0623'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '5.38'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.0994'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '97.59'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.451.98'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '4.07'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '15.59'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '7.32'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '3.00'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.279.46'

# Volume(1h) % change column (E) updates
$ws.Range("E3").Value = '  +1.14%  '
$ws.Range("E4").Value = '  -0.03%  '
$ws.Range("E5").Value = '  -0.68%  '
$ws.Range("E6").Value = '  -0.12%  '
$ws.Range("E7").Value = '  -0.02%  '
$ws.Range("E8").Value = '  +1.16%  '
$ws.Range("E9").Value = '  +1.60%  '
$ws.Range("E10").Value = '  +2.16%  '
$ws.Range("E11").Value = '  +2.94%  '
$ws.Range("E12").Value = '  +0.65%  '
$ws.Range("E13").Value = '  -1.42%  '
$ws.Range("E14").Value = '  +2.23%  '
$ws.Range("E15").Value = '  -1.03%  '
$ws.Range("E16").Value = '  +1.99%  '
$ws.Range("E17").Value = '  +1.59%  '
$ws.Range("E18").Value = '  +1.08%  '
$ws.Range("E19").Value = '  -3.09%  '
$ws.Range("E21").Value = '  +1.40%  '
$ws.Range("E22").Value = '  +0.76%  '
$ws.Range("E23").Value = '  -0.01%  '
$ws.Range("E24").Value = '  -1.20%  '
$ws.Range("E25").Value = '  -0.44%  '
$ws.Range("E26").Value = '  +1.08%  '
$ws.Range("E27").Value = '  +10.49%  '
$ws.Range("E28").Value = '  +1.80%  '
$ws.Range("E29").Value = '  -1.28%  '
$ws.Range("E30").Value = '  +2.34%  '
$ws.Range("E31").Value = '  +0.71%  '
$ws.Range("E32").Value = '  +3.78%  '
$ws.Range("E33").Value = '  +1.41%  '
$ws.Range("E34").Value = '  +0.81%  '
$ws.Range("E35").Value = '  +0.47%  '
$ws.Range("E36").Value = '  +4.24%  '
$ws.Range("E37").Value = '  +4.71%  '
$ws.Range("E38").Value = '  +0.11%  '
$ws.Range("E39").Value = '  -5.23%  '
$ws.Range("E40").Value = '  +6.65%  '
$ws.Range("E41").Value = '  -0.29%  '
$ws.Range("E42").Value = '  +1.65%  '
$ws.Range("E43").Value = '  +0.55%  '
$ws.Range("E44").Value = '  -0.70%  '
$ws.Range("E45").Value = '  -0.30%  '
$ws.Range("E46").Value = '  +3.56%  '
$ws.Range("E47").Value = '  -3.93%  '
$ws.Range("E48").Value = '  +3.59%  '
$ws.Range("E49").Value = '  +2.45%  '
$ws.Range("E50").Value = '  +1.83%  '
$ws.Range("E51").Value = '  +1.17%  '
